$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '28.194.93'
Set-TextCell 'E2' '  +0.07%  '
Set-TextCell 'D3' '1.870.47'
Set-TextCell 'E3' '  +2.01%  '
Set-TextCell 'D4' '1.001'
Set-TextCell 'E4' '  +0.25%  '
Set-TextCell 'D5' '311.62'
Set-TextCell 'E5' '  +0.08%  '
Set-TextCell 'E6' '  +0.08%  '
Set-TextCell 'D7' '0.5044'
Set-TextCell 'E7' '  -1.67%  '
Set-TextCell 'D8' '0.3918'
Set-TextCell 'E8' '  -0.99%  '
Set-TextCell 'D9' '0.09649'
Set-TextCell 'E9' '  -0.10%  '
Set-TextCell 'D10' '1.139'
Set-TextCell 'E10' '  +2.66%  '
Set-TextCell 'D11' '40.89'
Set-TextCell 'E11' '  -0.22%  '
Set-TextCell 'D12' '6.496'
Set-TextCell 'E12' '  +0.25%  '
Set-TextCell 'D13' '20.94'
Set-TextCell 'E13' '  +1.65%  '
Set-TextCell 'D14' '1.889.40'
Set-TextCell 'E14' '  +4.31%  '
Set-TextCell 'B15' 'Chainlink'
Set-TextCell 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D15' '7.430'
Set-TextCell 'E15' '  +0.37%  '
Set-TextCell 'B16' 'BinanceUSD'
Set-TextCell 'C16' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D16' '1.001'
Set-TextCell 'E16' '  +0.20%  '
Set-TextCell 'D17' '0.00001126'
Set-TextCell 'E17' '  -1.23%  '
Set-TextCell 'D18' '92.89'
Set-TextCell 'E18' '  -0.19%  '
Set-TextCell 'D19' '0.06628'
Set-TextCell 'E19' '  +0.33%  '
Set-TextCell 'E20' '  +0.72%  '
Set-TextCell 'E21' '  +0.13%  '
Set-TextCell 'D22' '6.147'
Set-TextCell 'E22' '  +1.30%  '
Set-TextCell 'D23' '28.261.29'
Set-TextCell 'E23' '  +0.11%  '
Set-TextCell 'D24' '11.34'
Set-TextCell 'E24' '  +1.53%  '
Set-TextCell 'D25' '2.281'
Set-TextCell 'E25' '  +1.26%  '
Set-TextCell 'D26' '2.530'
Set-TextCell 'E26' '  +3.28%  '
Set-TextCell 'D27' '2.091.55'
Set-TextCell 'E27' '  +2.97%  '
Set-TextCell 'D28' '21.20'
Set-TextCell 'E28' '  +2.10%  '
Set-TextCell 'D29' '157.30'
Set-TextCell 'E29' '  +0.02%  '
Set-TextCell 'D30' '127.22'
Set-TextCell 'E30' '  -1.58%  '
Set-TextCell 'B31' 'Stellar'
Set-TextCell 'C31' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D31' '0.1057'
Set-TextCell 'E31' '  -3.43%  '
Set-TextCell 'B32' 'ImmutableX'
Set-TextCell 'C32' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D32' '1.067'
Set-TextCell 'E32' '  +0.56%  '
Set-TextCell 'D33' '5.622'
Set-TextCell 'E33' '  -0.74%  '
Set-TextCell 'D34' '3.624'
Set-TextCell 'E34' '  -0.40%  '
Set-TextCell 'D35' '9.567'
Set-TextCell 'E35' '  +5.09%  '
Set-TextCell 'D36' '0.06752'
Set-TextCell 'E36' '  -2.75%  '
Set-TextCell 'E37' '  +1.64%  '
Set-TextCell 'E38' '  -0.22%  '
Set-TextCell 'E39' '  +1.30%  '
Set-TextCell 'D40' '11.46'
Set-TextCell 'E40' '  -1.35%  '
Set-TextCell 'D41' '4.982'
Set-TextCell 'E41' '  -1.00%  '
Set-TextCell 'D42' '1.178'
Set-TextCell 'E42' '  +1.87%  '
Set-TextCell 'E43' '  +0.16%  '
Set-TextCell 'E44' '  +1.77%  '
Set-TextCell 'D45' '0.6019'
Set-TextCell 'E45' '  +0.33%  '
Set-TextCell 'D46' '3.663'
Set-TextCell 'E46' '  -1.30%  '
Set-TextCell 'E47' '  -2.15%  '
Set-TextCell 'D48' '124.60'
Set-TextCell 'E48' '  -0.76%  '
Set-TextCell 'E49' '  +0.64%  '
Set-TextCell 'D50' '1.196'
Set-TextCell 'E50' '  +0.52%  '
Set-TextCell 'D51' '0.06834'
Set-TextCell 'E51' '  +0.63%  '
